$wb = $excel.ActiveWorkbook

# --- 1) Swap tab order: "review_info" becomes the first sheet, "hotel_info" second ---
$wsReview = $wb.Worksheets.Item("review_info")
$wsReview.Move($wb.Worksheets.Item(1))

# --- 2) Insert a new "State" column into hotel_info, right after "Hotel_Name" ---
# (re-fetch the sheet reference by name AFTER the move -- Move() re-targets
#  existing object handles by slot, so a reference captured before the move
#  would now point at the wrong sheet)
$wsHotel = $wb.Worksheets.Item("hotel_info")
$wsHotel.Columns.Item(3).Insert()
$wsHotel.Range("C1").Value = "State"
$wsHotel.Range("C2").Value = "Louisiana"
